$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.487.40'
$ws.Range("E2").Value = '  -5.11%  '
$ws.Range("D3").Value = '3.358.93'
$ws.Range("E3").Value = '  -3.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.36'
$ws.Range("E5").Value = '  -8.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '630.75'
$ws.Range("E6").Value = '  -6.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.37'
$ws.Range("E7").Value = '  -9.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.392'
$ws.Range("E8").Value = '  -9.28%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.942'
$ws.Range("E10").Value = '  -11.31%  '
$ws.Range("D11").Value = '3.356.67'
$ws.Range("E11").Value = '  -3.32%  '
$ws.Range("E12").Value = '  -7.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.48'
$ws.Range("E13").Value = '  -13.14%  '
$ws.Range("E14").Value = '  -2.97%  '
$ws.Range("D15").Value = '93.231.58'
$ws.Range("E15").Value = '  -5.25%  '
$ws.Range("D16").Value = '3.982.66'
$ws.Range("E16").Value = '  -3.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000244'
$ws.Range("E17").Value = '  -6.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.01'
$ws.Range("E18").Value = '  -12.33%  '
$ws.Range("D19").Value = '3.357.61'
$ws.Range("E19").Value = '  -3.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.90'
$ws.Range("E20").Value = '  -9.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.92'
$ws.Range("E21").Value = '  -6.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.00'
$ws.Range("E22").Value = '  -4.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.452'
$ws.Range("E23").Value = '  -16.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.14'
$ws.Range("E24").Value = '  -9.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000186'
$ws.Range("E25").Value = '  -8.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.31'
$ws.Range("E26").Value = '  -8.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.84'
$ws.Range("E27").Value = '  -8.59%  '
$ws.Range("D28").Value = '3.545.04'
$ws.Range("E28").Value = '  -3.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.50'
$ws.Range("E29").Value = '  -9.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.35'
$ws.Range("E30").Value = '  -8.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.65'
$ws.Range("E32").Value = '  -8.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.132'
$ws.Range("E33").Value = '  -9.48%  '
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.172'
$ws.Range("E35").Value = '  -10.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.70'
$ws.Range("E36").Value = '  -4.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.528'
$ws.Range("E37").Value = '  -9.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.50'
$ws.Range("E38").Value = '  -7.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '523.46'
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.40'
$ws.Range("E41").Value = '  -8.71%  '
$ws.Range("E42").Value = '  -5.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.879'
$ws.Range("E43").Value = '  -1.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.01'
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.67'
$ws.Range("E45").Value = '  -6.81%  '
$ws.Range("B46").Value = 'MantraDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.60'
$ws.Range("E46").Value = '  -1.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.49'
$ws.Range("E47").Value = '  -4.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.16'
$ws.Range("E48").Value = '  -3.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0395'
$ws.Range("E49").Value = '  -10.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.33'
$ws.Range("E50").Value = '  -6.72%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.13'
$ws.Range("E51").Value = '  -3.67%  '
